$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Right after "Write a minimum of 2 sentences ..." there is an empty
#    paragraph. Insert a brand-new list paragraph *before* it that answers
#    the prompt in red text (ListParagraph style, ilvl=1, numId=3 - same
#    list used by the surrounding bullets).
# ---------------------------------------------------------------------------
$wOpenXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$targetEmpty = $d.Content.Find
$targetEmpty.Execute("Write a minimum of 2 sentences") | Out-Null

$i = 0
$afterWriteIdx = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Write a minimum of 2 sentences*") {
        $afterWriteIdx = $i + 1
    }
}

$pEmptyAfterWrite = $d.Paragraphs($afterWriteIdx)
$pEmptyAfterWrite.Range.InsertParagraphBefore() | Out-Null
$pAns = $d.Paragraphs($afterWriteIdx)

$ansXml = '<w:p ' + $wOpenXmlNs + '>' `
    + '<w:pPr>' `
        + '<w:pStyle w:val="ListParagraph"/>' `
        + '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>' `
        + '<w:tabs><w:tab w:val="left" w:pos="4176"/></w:tabs>' `
    + '</w:pPr>' `
    + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr>' `
        + '<w:t>Ans: Most people in this group have normal blood pressure</w:t>' `
    + '</w:r>' `
    + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr>' `
        + '<w:t>, followed by hypertension stage 1.  Only 1 person has hypertensive crisis. (There are plenty of other things to write.)</w:t>' `
    + '</w:r>' `
    + '</w:p>'
$pAns.Range.InsertXML($ansXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) After the "Have students work on their own, ..." paragraph, insert a
#    blank paragraph followed by a new paragraph that tells TAs how to
#    verify submissions (SpeedGrader note).
# ---------------------------------------------------------------------------
$j = 0
$haveStudentsIdx = -1
foreach ($p in $d.Paragraphs) {
    $j = $j + 1
    if ($p.Range.Text -like "Have students work on their own*") {
        $haveStudentsIdx = $j
    }
}

$pHave = $d.Paragraphs($haveStudentsIdx)
$pHave.Range.InsertParagraphAfter() | Out-Null
$pBlank = $d.Paragraphs($haveStudentsIdx + 1)

$blankXml = '<w:p ' + $wOpenXmlNs + '>' `
    + '<w:pPr><w:tabs><w:tab w:val="left" w:pos="4176"/></w:tabs></w:pPr>' `
    + '</w:p>'
$pBlank.Range.InsertXML($blankXml) | Out-Null

$pBlank2 = $d.Paragraphs($haveStudentsIdx + 1)
$pBlank2.Range.InsertParagraphAfter() | Out-Null
$pSpeed = $d.Paragraphs($haveStudentsIdx + 2)

$speedXml = '<w:p ' + $wOpenXmlNs + '>' `
    + '<w:pPr><w:tabs><w:tab w:val="left" w:pos="4176"/></w:tabs></w:pPr>' `
    + '<w:r><w:lastRenderedPageBreak/>' `
        + '<w:t xml:space="preserve">You can verify their submission by going to the Grades, or </w:t>' `
    + '</w:r>' `
    + '<w:proofErr w:type="spellStart"/>' `
    + '<w:r><w:t>SpeedGrader</w:t></w:r>' `
    + '<w:proofErr w:type="spellEnd"/>' `
    + '<w:r><w:t>.</w:t></w:r>' `
    + '</w:p>'
$pSpeed.Range.InsertXML($speedXml) | Out-Null

Write-Output "Edit applied."
